# Update achievement. Add links.
# Fills in newly-recorded grade values across several students' rows
# (columns I/J/L/H for assignments A6/A7/A9/A5, W/X/V for КР3/КР4/КР2, AH for Л1),
# then moves the active selection to K35 (matching the saved workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 0.1

$ws.Range("I3").Value = 0.1
$ws.Range("X3").Value = 0.05

$ws.Range("X4").Value = 0.1

$ws.Range("I5").Value = 0.6

$ws.Range("J7").Value = 0.6
$ws.Range("X7").Value = 0.1

$ws.Range("X8").Value = 0.1

$ws.Range("X9").Value = 0.1

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0.6
$ws.Range("X10").Value = 0.1

$ws.Range("I11").Value = 0.05

$ws.Range("I12").Value = 0.5
$ws.Range("X12").Value = 0.1

$ws.Range("X13").Value = 0.1

$ws.Range("X14").Value = 0.1

$ws.Range("H15").Value = 0.4
$ws.Range("J15").Value = 0.3

$ws.Range("I16").Value = 0.8
$ws.Range("J16").Value = 0.6
$ws.Range("X16").Value = 0.1

$ws.Range("I17").Value = 0.6
$ws.Range("X17").Value = 0.1

$ws.Range("I19").Value = 1
$ws.Range("X19").Value = 0.1

$ws.Range("J20").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("W20").Value = 0.1
$ws.Range("X20").Value = 0.1
$ws.Range("AH20").Value = 1

$ws.Range("I22").Value = 1
$ws.Range("X22").Value = 0.1

$ws.Range("J24").Value = 1
$ws.Range("X24").Value = 0.1

$ws.Range("I26").Value = 1.5
$ws.Range("J26").Value = 1
$ws.Range("X26").Value = 0.1

$ws.Range("J27").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("W27").Value = 0.1
$ws.Range("X27").Value = 0.2

$ws.Range("J28").Value = 0.6
$ws.Range("X28").Value = 0.1

$ws.Range("I29").Value = 0.6
$ws.Range("J29").Value = 0.1

$ws.Range("J31").Value = 0.6
$ws.Range("X31").Value = 0.2

$ws.Range("I32").Value = 0.3
$ws.Range("J32").Value = 0.1

$ws.Range("J34").Value = 1
$ws.Range("X34").Value = 0.1

$ws.Range("J35").Value = 1
$ws.Range("V35").Value = 0.2
$ws.Range("W35").Value = 0.2
$ws.Range("X35").Value = 0.2

$ws.Range("K35").Select()
